# Recreate the "added results for logistic regression" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width: 13.5 -> 16.5 (character units). The ColumnWidth COM
# property uses a slightly different unit than the raw OOXML width, so use
# a calibrated input that yields the target stored width.
$ws.Columns.Item(1).ColumnWidth = 15.65

# --- New section header label in A24: "Logistic Regression"
$ws.Range("A24").Value = "Logistic Regression"

# --- Row 25: duplicate the header row (row 2) formats + values exactly
$ws.Range("A2:K2").Copy($ws.Range("A25"))

# --- Rows 26-41: new Logistic Regression weekly data (plain values, no
# special styling - matches the source diff which has no "s" attribute).
$logisticRows = @(
  @(2,3,0,4,3,3,3,3,3,3,3),
  @(3,4,3,5,4,4,4,4,4,4,4),
  @(4,2,2,1,2,3,2,3,3,3,3),
  @(5,3,2,3,3,4,4,4,4,4,4),
  @(6,3,0,1,3,4,2,2,3,3,3),
  @(7,0,2,2,0,0,2,2,2,2,1),
  @(8,2,1,1,2,3,3,4,3,3,4),
  @(9,2,4,3,2,2,2,3,3,3,3),
  @(10,2,3,1,2,2,2,2,2,2,2),
  @(11,3,3,3,3,3,4,4,4,4,4),
  @(12,2,3,1,2,1,2,2,2,2,0),
  @(13,2,2,2,2,2,3,2,2,2,2),
  @(14,2,1,3,2,3,2,2,2,2,2),
  @(15,1,0,3,1,3,2,2,2,2,2),
  @(16,2,2,2,2,2,3,4,4,4,3),
  @(17,3,0,3,3,2,4,4,4,4,4)
)

$r = 26
foreach ($rowVals in $logisticRows) {
    $c = 1
    foreach ($v in $rowVals) {
        $ws.Cells.Item($r, $c).Value = $v
        $c = $c + 1
    }
    $r = $r + 1
}

# --- Row 43: totals row (label + SUM formulas over the new data block)
$ws.Range("A43").Value = "Total Wins:"
$cols = @("B","C","D","E","F","G","H","I","J","K")
foreach ($col in $cols) {
    $ws.Range($col + "43").Formula = "=SUM(" + $col + "26:" + $col + "41)"
}

# --- Selection / scroll position bookkeeping, best effort.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("K44").Select()
